$d = $word.ActiveDocument

# Delete everything after the first paragraph ("[!] Needs Review"),
# i.e. remove the Q1/Q4/Q5 question & answer paragraphs while keeping
# the heading paragraph and the final section properties.
$firstPara = $d.Paragraphs(1)
$startDelete = $firstPara.Range.End
$endDelete = $d.Content.End

$deleteRange = $d.Range($startDelete, $endDelete)
$deleteRange.Delete()
